# Adds two new tracked datasets ("GFS Social Expenditure" in column AG,
# "Tax Wedge" in column AH) to the "Datasets and Years" country-year index.
#
# Column AG gets a checkmark for every data row (2000-2017 originally had
# rows 2..29 -> years 1990..2017), column AH only gets a checkmark for the
# more recent rows (2000-2017, i.e. rows 12..29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkmark = [char]0x2713   # "✓" (matches the shared string already used for s="12")

# --- Copy the existing header/body style (s="1") onto the new cells first,
#     so the values we set below land with the same formatting as the rest
#     of the sheet instead of picking up a brand-new default style.
[void]$ws.Range("A1").Copy()
[void]$ws.Range("AG1:AG29").PasteSpecial(-4122)   # xlPasteFormats

[void]$ws.Range("A1").Copy()
[void]$ws.Range("AH1").PasteSpecial(-4122)        # xlPasteFormats

[void]$ws.Range("A1").Copy()
[void]$ws.Range("AH12:AH29").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# --- Headers
$ws.Range("AG1").Value = "GFS Social Expenditure"
$ws.Range("AH1").Value = "Tax Wedge"

# --- Column AG: checkmark for every year row (2..29)
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 33).Value = $checkmark   # column 33 = AG
}

# --- Column AH: checkmark only for rows 12..29 (years 2000-2017)
for ($r = 12; $r -le 29; $r++) {
    $ws.Cells.Item($r, 34).Value = $checkmark   # column 34 = AH
}

# --- Pane / selection bookkeeping to mirror the author's final view state.
$ws.Activate()
[void]$ws.Range("AH12:AH29").Select()
